$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 8 de Julio de 2020 a las 11:28"

# Update country data rows (refreshed stats + reordering of a few neighboring countries)
# Row 4
$ws.Range("B4").Value2 = 3097421
$ws.Range("C4").Value2 = 337
$ws.Range("E4").Value2 = 1608363

# Row 6
$ws.Range("B6").Value2 = 746139
$ws.Range("C6").Value2 = 2658
$ws.Range("D6").Value2 = 458489
$ws.Range("E6").Value2 = 266967
$ws.Range("G6").Value2 = 30
$ws.Range("H6").Value2 = 20683

# Row 20
$ws.Range("A20").Value2 = "Banglades"
$ws.Range("B20").Value2 = 172134
$ws.Range("C20").Value2 = 3489
$ws.Range("D20").Value2 = 80838
$ws.Range("E20").Value2 = 89099
$ws.Range("G20").Value2 = 46
$ws.Range("H20").Value2 = 2197

# Row 21
$ws.Range("A21").Value2 = "Francia"
$ws.Range("B21").Value2 = 168810
$ws.Range("D21").Value2 = 77655
$ws.Range("E21").Value2 = 61222
$ws.Range("H21").Value2 = 29933

# Row 29
$ws.Range("B29").Value2 = 68079
$ws.Range("C29").Value2 = 1853
$ws.Range("D29").Value2 = 31585
$ws.Range("E29").Value2 = 33135
$ws.Range("G29").Value2 = 50
$ws.Range("H29").Value2 = 3359

# Row 39
$ws.Range("B39").Value2 = 50207
$ws.Range("C39").Value2 = 1210
$ws.Range("D39").Value2 = 32005
$ws.Range("E39").Value2 = 17969
$ws.Range("G39").Value2 = 9
$ws.Range("H39").Value2 = 233

# Row 46
$ws.Range("B46").Value2 = 36689
$ws.Range("C46").Value2 = 277
$ws.Range("D46").Value2 = 24878
$ws.Range("E46").Value2 = 10269
$ws.Range("G46").Value2 = 14
$ws.Range("H46").Value2 = 1542

# Row 48
$ws.Range("A48").Value2 = "Israel"
$ws.Range("B48").Value2 = 32714
$ws.Range("C48").Value2 = 492
$ws.Range("D48").Value2 = 18267
$ws.Range("E48").Value2 = 14104
$ws.Range("G48").Value2 = 1
$ws.Range("H48").Value2 = 343

# Row 49
$ws.Range("A49").Value2 = "Suiza"
$ws.Range("B49").Value2 = 32369
$ws.Range("D49").Value2 = 29300
$ws.Range("E49").Value2 = 1103
$ws.Range("H49").Value2 = 1966

# Row 66
$ws.Range("B66").Value2 = 14730
$ws.Range("C66").Value2 = 123
$ws.Range("D66").Value2 = 10848
$ws.Range("E66").Value2 = 3642

# Row 82
$ws.Range("B82").Value2 = 7265
$ws.Range("C82").Value2 = 3
$ws.Range("E82").Value2 = 236

# Row 93
$ws.Range("A93").Value2 = "Estado de Palestina"
$ws.Range("B93").Value2 = 5029
$ws.Range("C93").Value2 = 382
$ws.Range("D93").Value2 = 494
$ws.Range("E93").Value2 = 4516
$ws.Range("G93").Value2 = 1
$ws.Range("H93").Value2 = 19

# Row 94
$ws.Range("A94").Value2 = "Mauritania"
$ws.Range("B94").Value2 = 5024
$ws.Range("D94").Value2 = 1944
$ws.Range("E94").Value2 = 2945
$ws.Range("H94").Value2 = 135

# Row 95
$ws.Range("A95").Value2 = "Republica de Yibuti"
$ws.Range("B95").Value2 = 4878
$ws.Range("D95").Value2 = 4621
$ws.Range("E95").Value2 = 202
$ws.Range("H95").Value2 = 55

# Row 122
$ws.Range("B122").Value2 = 1763
$ws.Range("C122").Value2 = 24
$ws.Range("D122").Value2 = 1429
$ws.Range("E122").Value2 = 223

# Row 127
$ws.Range("B127").Value2 = 1324
$ws.Range("C127").Value2 = 24
$ws.Range("D127").Value2 = 1167
$ws.Range("E127").Value2 = 150

# Row 134
$ws.Range("B134").Value2 = 1141
$ws.Range("C134").Value2 = 7
$ws.Range("E134").Value2 = 103

# Row 140
$ws.Range("B140").Value2 = 977
$ws.Range("C140").Value2 = 6
$ws.Range("D140").Value2 = 904
$ws.Range("E140").Value2 = 73

# Row 155
$ws.Range("A155").Value2 = "Namibia"
$ws.Range("B155").Value2 = 593
$ws.Range("C155").Value2 = 54
$ws.Range("D155").Value2 = 25
$ws.Range("E155").Value2 = 568
$ws.Range("H155").Value2 = 0

# Row 156
$ws.Range("A156").Value2 = "Reunion"
$ws.Range("B156").Value2 = 551
$ws.Range("D156").Value2 = 472
$ws.Range("E156").Value2 = 76
$ws.Range("H156").Value2 = 3

# Row 209
$ws.Range("A209").Value2 = "Groenlandia"

# Row 210
$ws.Range("A210").Value2 = "Islas Malvinas"
